$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1661.0714
$ws.Range("J28").Value = 5248.8335
$ws.Range("L28").Value = 5248.8335
$ws.Range("N28").Value = -6218.8335
$ws.Range("H45").Value = 8014.2856
$ws.Range("J45").Value = 9183.333000000001
$ws.Range("L45").Value = 27549.999
$ws.Range("N45").Value = -27933.999
$ws.Range("H62").Value = 12084
$ws.Range("I62").Value = 12861.5
$ws.Range("K62").Value = 12861.5
$ws.Range("M62").Value = -12237.5
$ws.Range("H65").Value = 12084
$ws.Range("I65").Value = 12861.5
$ws.Range("K65").Value = 64307.5
$ws.Range("M65").Value = -61187.5
$ws.Range("H107").Value = 735.55554
$ws.Range("I107").Value = 702.625
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 702.625
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1217.375
$ws.Range("N107").Value = -4839
$ws.Range("H135").Value = 498.8421
$ws.Range("I135").Value = 483.6
$ws.Range("K135").Value = 4352.400000000001
$ws.Range("M135").Value = -1817.400000000001
$ws.Range("H137").Value = 5746.3794
$ws.Range("I137").Value = 1849.1765
$ws.Range("J137").Value = 11267.417
$ws.Range("K137").Value = 5547.529500000001
$ws.Range("L137").Value = 33802.251
$ws.Range("M137").Value = -2997.529500000001
$ws.Range("N137").Value = -38902.251
$ws.Range("H138").Value = 394856.3
$ws.Range("J138").Value = 514626.34
$ws.Range("L138").Value = 1543879.02
$ws.Range("N138").Value = -1554159.02
$ws.Range("H141").Value = 745.75
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1211.2106
$ws.Range("I2").Value = 1144.5625
$ws.Range("K2").Value = 1144.5625
$ws.Range("M2").Value = -1031.5625
$ws.Range("H16").Value = 4356.4443
$ws.Range("I16").Value = 1399.6666
$ws.Range("J16").Value = 5834.8335
$ws.Range("K16").Value = 1399.6666
$ws.Range("L16").Value = 5834.8335
$ws.Range("M16").Value = -1112.6666
$ws.Range("N16").Value = -6408.8335
$ws.Range("H32").Value = 5277.7
$ws.Range("I32").Value = 4722.154
$ws.Range("K32").Value = 4722.154
$ws.Range("M32").Value = -4435.154
$ws.Range("H61").Value = 5268.5483
$ws.Range("I61").Value = 2628.7646
$ws.Range("K61").Value = 2628.7646
$ws.Range("M61").Value = -2416.7646
$ws.Range("H116").Value = 1211.2106
$ws.Range("I116").Value = 1144.5625
$ws.Range("K116").Value = 1144.5625
$ws.Range("M116").Value = 1149.4375
$ws.Range("H122").Value = 3240.5312
$ws.Range("I122").Value = 3023.8928
$ws.Range("K122").Value = 9071.678400000001
$ws.Range("M122").Value = -6621.678400000001
$ws.Range("H136").Value = 5268.5483
$ws.Range("I136").Value = 2628.7646
$ws.Range("K136").Value = 7886.293799999999
$ws.Range("M136").Value = -5336.293799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1211.2106
$ws.Range("I3").Value = 1144.5625
$ws.Range("K3").Value = 1144.5625
$ws.Range("M3").Value = -1030.5625
$ws.Range("H80").Value = 615.0714
$ws.Range("I80").Value = 1217.5
$ws.Range("J80").Value = 374.1
$ws.Range("K80").Value = 1217.5
$ws.Range("L80").Value = 374.1
$ws.Range("M80").Value = -219.5
$ws.Range("N80").Value = -2370.1
$ws.Range("H83").Value = 615.0714
$ws.Range("I83").Value = 1217.5
$ws.Range("J83").Value = 374.1
$ws.Range("K83").Value = 6087.5
$ws.Range("L83").Value = 1870.5
$ws.Range("M83").Value = -1095.5
$ws.Range("N83").Value = -11854.5
$ws.Range("H134").Value = 2862.6667
$ws.Range("I134").Value = 2405.1667
$ws.Range("K134").Value = 7215.500100000001
$ws.Range("M134").Value = -4680.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4028.6724
$ws.Range("I31").Value = 3365.3713
$ws.Range("J31").Value = 5038.0435
$ws.Range("K31").Value = 3365.3713
$ws.Range("L31").Value = 5038.0435
$ws.Range("M31").Value = -3070.3713
$ws.Range("N31").Value = -5628.0435
$ws.Range("H34").Value = 4028.6724
$ws.Range("I34").Value = 3365.3713
$ws.Range("J34").Value = 5038.0435
$ws.Range("K34").Value = 3365.3713
$ws.Range("L34").Value = 5038.0435
$ws.Range("M34").Value = -3163.3713
$ws.Range("N34").Value = -5442.0435
$ws.Range("H58").Value = 2672.111
$ws.Range("I58").Value = 1532
$ws.Range("K58").Value = 1532
$ws.Range("M58").Value = -1329
$ws.Range("H99").Value = 6329.7
$ws.Range("I99").Value = 6908.696
$ws.Range("J99").Value = 4427.2856
$ws.Range("K99").Value = 6908.696
$ws.Range("L99").Value = 4427.2856
$ws.Range("M99").Value = -5410.696
$ws.Range("N99").Value = -7423.2856
$ws.Range("H107").Value = 635.25
$ws.Range("I107").Value = 700
$ws.Range("K107").Value = 700
$ws.Range("M107").Value = 1220
$ws.Range("H126").Value = 6329.7
$ws.Range("I126").Value = 6908.696
$ws.Range("J126").Value = 4427.2856
$ws.Range("K126").Value = 20726.088
$ws.Range("L126").Value = 13281.8568
$ws.Range("M126").Value = -18256.088
$ws.Range("N126").Value = -18221.8568
$ws.Range("H132").Value = 5756.625
$ws.Range("I132").Value = 5305.222
$ws.Range("K132").Value = 15915.666
$ws.Range("M132").Value = -13385.666
$ws.Range("H134").Value = 2409.1333
$ws.Range("I134").Value = 2094.4
$ws.Range("J134").Value = 3982.8
$ws.Range("K134").Value = 6283.200000000001
$ws.Range("L134").Value = 11948.4
$ws.Range("M134").Value = -3748.200000000001
$ws.Range("N134").Value = -17018.4
$ws.Range("H136").Value = 2672.111
$ws.Range("I136").Value = 1532
$ws.Range("K136").Value = 4596
$ws.Range("M136").Value = -2046
$ws.Range("H141").Value = 299608.94
$ws.Range("J141").Value = 299608.94
$ws.Range("L141").Value = 299608.94
$ws.Range("N141").Value = -309968.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3114.5
$ws.Range("I5").Value = 476.4
$ws.Range("K5").Value = 1429.2
$ws.Range("M5").Value = -1317.2
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 6000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -6554
$ws.Range("H69").Value = 2489.2856
$ws.Range("J69").Value = 2829.1667
$ws.Range("L69").Value = 8487.500100000001
$ws.Range("N69").Value = -10109.5001
$ws.Range("H72").Value = 2489.2856
$ws.Range("J72").Value = 2829.1667
$ws.Range("L72").Value = 25462.5003
$ws.Range("N72").Value = -33574.5003
$ws.Range("H122").Value = 1173
$ws.Range("I122").Value = 1139.5
$ws.Range("J122").Value = 1178.36
$ws.Range("K122").Value = 10255.5
$ws.Range("L122").Value = 10605.24
$ws.Range("M122").Value = -7805.5
$ws.Range("N122").Value = -15505.24
$ws.Range("H131").Value = 2435.5557
$ws.Range("J131").Value = 2050
$ws.Range("L131").Value = 6150
$ws.Range("N131").Value = -16230
$ws.Range("H132").Value = 3006.9656
$ws.Range("I132").Value = 1978.7142
$ws.Range("J132").Value = 3966.6667
$ws.Range("K132").Value = 17808.4278
$ws.Range("L132").Value = 35700.0003
$ws.Range("M132").Value = -15278.4278
$ws.Range("N132").Value = -40760.0003
$ws.Range("H135").Value = 3114.5
$ws.Range("I135").Value = 476.4
$ws.Range("K135").Value = 4287.599999999999
$ws.Range("M135").Value = -1752.599999999999
$ws.Range("H141").Value = 18198
$ws.Range("I141").Value = 9496.666999999999
$ws.Range("K141").Value = 28490.001
$ws.Range("M141").Value = -23310.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4945.7915
$ws.Range("I113").Value = 5761.8125
$ws.Range("K113").Value = 5761.8125
$ws.Range("M113").Value = -3591.8125
$ws.Range("H126").Value = 3675.3125
$ws.Range("I126").Value = 2775.182
$ws.Range("J126").Value = 5655.6
$ws.Range("K126").Value = 8325.545999999998
$ws.Range("L126").Value = 16966.8
$ws.Range("M126").Value = -5855.545999999998
$ws.Range("N126").Value = -21906.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1553.45
$ws.Range("J46").Value = 1460.9032
$ws.Range("L46").Value = 1460.9032
$ws.Range("N46").Value = -1836.9032
$ws.Range("H100").Value = 4851.391
$ws.Range("I100").Value = 4455.6665
$ws.Range("K100").Value = 4455.6665
$ws.Range("M100").Value = -3914.6665
$ws.Range("H122").Value = 5119.5
$ws.Range("I122").Value = 5082.5454
$ws.Range("K122").Value = 15247.6362
$ws.Range("M122").Value = -12797.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1120.1904
$ws.Range("I113").Value = 1151.9375
$ws.Range("K113").Value = 3455.8125
$ws.Range("M113").Value = -1285.8125
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080
$ws.Range("H136").Value = 38465640
$ws.Range("I136").Value = 47620790
$ws.Range("J136").Value = 13997.4
$ws.Range("K136").Value = 142862370
$ws.Range("L136").Value = 41992.2
$ws.Range("M136").Value = -142859820
$ws.Range("N136").Value = -47092.2
